# Edit Collaboration.docx:
#  1. Remove the <w:tblPrEx> (row cell-margin override) block from both
#     table rows.
#  2. Split the "Game Screen design elements" run into two runs so the
#     text reads "Some Game Screen design elements" (new leading run
#     holds "Some " with a trailing space, xml:space="preserve", matching
#     the formatting of the original run).
#
# Word's object model has no direct row-level "clear tblPrEx" verb, and
# Range.WordOpenXML is read-only (confirmed by this runtime), so we pull
# the document's full OOXML, edit it as text, and push it back with
# Range.InsertXML on the whole-document Content range -- the one range
# InsertXML treats as a full replace rather than an insert/paste.

$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

# --- 1. Strip the per-row tblPrEx overrides (top/bottom cell margin = 0) ---
$tblPrEx = '<w:tblPrEx><w:tblCellMar><w:top w:w="0" w:type="dxa"/><w:bottom w:w="0" w:type="dxa"/></w:tblCellMar></w:tblPrEx>'
$tblPrExCount = ([regex]::Matches($xml, [regex]::Escape($tblPrEx))).Count
if ($tblPrExCount -ne 2) {
    throw "expected 2 tblPrEx blocks, found $tblPrExCount"
}
$xml = $xml.Replace($tblPrEx, "")

# --- 2. Split "Game Screen design elements" run, prefixing "Some " ---
$oldRun = '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/></w:rPr><w:t>Game Screen design elements</w:t></w:r>'
$oldRunCount = ([regex]::Matches($xml, [regex]::Escape($oldRun))).Count
if ($oldRunCount -ne 1) {
    throw "expected 1 'Game Screen design elements' run, found $oldRunCount"
}
$newRun = '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Some </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/></w:rPr><w:t>Game Screen design elements</w:t></w:r>'
$xml = $xml.Replace($oldRun, $newRun)

$d.Content.InsertXML($xml) | Out-Null
